$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C4"  = -12.44299999999999
    "B9"  = 5.418699999999996
    "C9"  = -10.82249999999999
    "C11" = -12.6218
    "B13" = 6.713599999999996
    "B16" = 4.925199999999996
    "C16" = -13.5826
    "B18" = 6.494299999999998
    "B20" = 8.856299999999994
    "C23" = -12.2571
    "C24" = -12.8202
    "B26" = 4.944300000000004
    "C26" = -13.10710000000001
    "B27" = 6.231400000000004
    "B29" = 5.343299999999999
    "C34" = -11.89320000000001
    "B35" = 8.669299999999994
    "C35" = -13.30580000000001
    "B36" = 9.35150000000001
    "C44" = -13.1925
    "B45" = 4.737500000000004
    "C48" = -11.34609999999999
    "C49" = -13.6017
    "C52" = -10.9172
    "B55" = 7.007699999999994
    "B57" = 4.979599999999994
    "C66" = -11.2329
    "C67" = -10.87719999999999
    "B69" = 5.424699999999993
    "C73" = -11.18040000000001
    "B76" = 4.921499999999998
    "B78" = 9.565599999999996
    "C78" = -12.26670000000001
    "C80" = -11.65640000000001
    "B82" = 5.2512
    "B83" = 5.1726
    "C91" = -12.925
    "B93" = 5.353400000000001
    "B97" = 5.511799999999997
    "C97" = -11.1068
    "C99" = -12.8511
    "C104" = -12.70200000000002
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
